$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'" + '69.422.70'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'" + '  +2.91%  '
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'" + '3.392.14'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'" + '  +2.20%  '
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'" + '  +0.07%  '
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'" + '588.71'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'" + '  +2.21%  '
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'" + '180.88'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'" + '  +3.95%  '
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'" + '  -0.04%  '
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'" + '  +1.19%  '
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'" + '  +10.67%  '
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'" + '0.589'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'" + '  +1.83%  '
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'" + '48.87'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'" + '  +6.80%  '
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'" + '  +5.50%  '
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'" + '684.49'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'" + '  -2.49%  '
$ws.Range("E13").Style = "Normal"
$ws.Range("B14").Value = "'" + 'Polkadot'
$ws.Range("B14").Style = "Normal"
$ws.Range("C14").Value = "'" + 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("C14").Style = "Normal"
$ws.Range("D14").Value = "'" + '8.61'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'" + '  +2.46%  '
$ws.Range("E14").Style = "Normal"
$ws.Range("B15").Value = "'" + 'WrappedliquidstakedEther2.0'
$ws.Range("B15").Style = "Normal"
$ws.Range("C15").Value = "'" + 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("C15").Style = "Normal"
$ws.Range("D15").Value = "'" + '3.940.17'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'" + '  +2.15%  '
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'" + '69.451.28'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'" + '  +2.93%  '
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'" + '3.399.12'
$ws.Range("D17").Style = "Normal"
$ws.Range("D19").Value = "'" + '17.72'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'" + '  +2.15%  '
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'" + '11.40'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'" + '  +4.21%  '
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'" + '  +1.68%  '
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'" + '  +0.57%  '
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'" + '  +1.77%  '
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'" + '104.65'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'" + '  +6.11%  '
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'" + '  +3.21%  '
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'" + '  +2.23%  '
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'" + '9.64'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'" + '  +3.59%  '
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'" + '34.60'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'" + '  +4.04%  '
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'" + '8.71'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'" + '  +2.66%  '
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'" + '  -0.59%  '
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'" + '11.19'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'" + '  +2.19%  '
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'" + '3.69'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'" + '  +11.92%  '
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'" + '555.25'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'" + '  -2.22%  '
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'" + '  +1.94%  '
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'" + '58.46'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'" + '  +2.93%  '
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'" + '3.725.65'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'" + '  +0.93%  '
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'" + '  +0.09%  '
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'" + '  +8.80%  '
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'" + '35.04'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'" + '  +2.63%  '
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'" + '0.0₃0714'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'" + '  +7.22%  '
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'" + '3.25'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'" + '  +3.26%  '
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'" + '  +2.74%  '
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'" + '  +2.37%  '
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'" + '  +3.91%  '
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'" + '  -0.75%  '
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'" + '  -0.73%  '
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'" + '  +1.60%  '
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'" + '1.40'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'" + '  +7.02%  '
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'" + '  +0.08%  '
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'" + '132.77'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'" + '  +2.99%  '
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'" + '2.64'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'" + '  -1.48%  '
$ws.Range("E51").Style = "Normal"
